$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so numeric-looking strings are not
# auto-converted to floating point numbers (values must stay verbatim text,
# matching the source data which uses inline strings for these columns).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.403.42'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '1.917.03'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '240.87'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '0.4675'
$ws.Range('E7').Value = '  -2.40%  '
$ws.Range('D8').Value = '0.2843'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').Value = '0.06829'
$ws.Range('E9').Value = '  +4.90%  '
$ws.Range('D10').Value = '107.31'
$ws.Range('E10').Value = '  +13.80%  '
$ws.Range('D11').Value = '17.96'
$ws.Range('E11').Value = '  -4.95%  '
$ws.Range('D12').Value = '1.908.31'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '0.07619'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '5.178'
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').Value = '0.6539'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '284.81'
$ws.Range('E16').Value = '  -4.38%  '
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.000007584'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '12.94'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').Value = '2.160.18'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = '5.202'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').Value = '6.190'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').Value = '168.05'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = '9.230'
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('D27').Value = '21.37'
$ws.Range('E27').Value = '  +9.18%  '
$ws.Range('D28').Value = '2.037'
$ws.Range('E28').Value = '  +4.50%  '
$ws.Range('D29').Value = '0.1069'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').Value = '1.369'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').Value = '4.129'
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').Value = '3.935'
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('D33').Value = '0.05030'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '0.7361'
$ws.Range('E34').Value = '  +1.86%  '
$ws.Range('D35').Value = '1.143'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').Value = '0.02018'
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('D40').Value = '2.043'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').Value = '0.8730'
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('D43').Value = '5.817'
$ws.Range('E43').Value = '  +4.11%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').Value = '52.77'
$ws.Range('E45').Value = '  +25.71%  '
$ws.Range('D46').Value = '0.4190'
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').Value = '67.32'
$ws.Range('E47').Value = '  +1.68%  '
$ws.Range('D48').Value = '7.138'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('D49').Value = '9.168'
$ws.Range('E49').Value = '  +3.13%  '
$ws.Range('D50').Value = '0.1205'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').Value = '34.54'
$ws.Range('E51').Value = '  -0.38%  '
